$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: D=Fecha(4), M=Volumen(13), N=Precio minimo(14), O=Precio maximo(15), P=Precio promedio ponderado(16), S=Precio $/Kg(19)
$data = @(
    @(2, 44336, 100, 19500, 20000, 19750, 988),
    @(3, 44365, 100, 20000, 21000, 20500, 1025),
    @(4, 44879, 100, 28000, 30000, 29000, 1450),
    @(5, 44343, 100, 19500, 20000, 19750, 988),
    @(6, 44410, 200, 20000, 21000, 20500, 1025),
    @(7, 44448, 100, 20000, 21000, 20500, 1025),
    @(8, 44874, 240, 29000, 30000, 29500, 1475),
    @(9, 44442, 140, 20000, 21000, 20500, 1025),
    @(10, 44427, 200, 20000, 21000, 20500, 1025),
    @(11, 44882, 120, 28000, 30000, 29000, 1450),
    @(12, 44473, 40, 19500, 20000, 19750, 988),
    @(13, 44326, 160, 19500, 20000, 19750, 988),
    @(14, 44809, 60, 27000, 28000, 27500, 1375),
    @(15, 44467, 200, 20000, 21000, 20500, 1025),
    @(16, 44407, 160, 20000, 21000, 20500, 1025),
    @(17, 44364, 140, 20000, 21000, 20500, 1025),
    @(18, 44435, 260, 20000, 22000, 21115, 1056),
    @(19, 44301, 100, 18000, 19000, 18500, 925),
    @(20, 44420, 160, 20000, 21000, 20500, 1025),
    @(21, 44441, 160, 20000, 21000, 20500, 1025),
    @(22, 44784, 160, 27000, 28000, 27500, 1375),
    @(23, 44462, 100, 19500, 20000, 19750, 988),
    @(24, 44350, 160, 19000, 20000, 19500, 975),
    @(25, 44474, 200, 19000, 20000, 19500, 975),
    @(26, 44418, 200, 20000, 21000, 20500, 1025),
    @(27, 44445, 160, 20000, 21000, 20500, 1025),
    @(28, 44417, 160, 20000, 21000, 20500, 1025),
    @(29, 44781, 160, 23000, 24000, 23500, 1175),
    @(30, 44434, 100, 20000, 21000, 20500, 1025),
    @(31, 44428, 100, 20000, 21000, 20500, 1025),
    @(32, 44880, 100, 28000, 30000, 29000, 1450),
    @(33, 44810, 100, 27000, 28000, 27500, 1375),
    @(34, 44466, 100, 20000, 21000, 20500, 1025),
    @(35, 44782, 200, 23500, 24000, 23750, 1188),
    @(36, 44776, 160, 23000, 24000, 23500, 1175),
    @(37, 44431, 160, 21000, 22000, 21500, 1075),
    @(38, 44315, 100, 20000, 21000, 20500, 1025),
    @(39, 44335, 200, 19000, 20000, 19500, 975),
    @(40, 44333, 100, 19500, 20000, 19750, 988),
    @(41, 44778, 100, 23000, 24000, 23500, 1175)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 13).Value = $row[2]
    $ws.Cells.Item($r, 14).Value = $row[3]
    $ws.Cells.Item($r, 15).Value = $row[4]
    $ws.Cells.Item($r, 16).Value = $row[5]
    $ws.Cells.Item($r, 19).Value = $row[6]
}

Write-Output "Updated $($data.Count) rows"